$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text (non-numeric-looking) cell updates: coin name / link columns
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'

# Volume/percentage column (E): plain text assignment is safe since values
# contain non-numeric characters ("%" and surrounding spaces).
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('E7').Value = '  -1.15%  '
$ws.Range('E8').Value = '  -2.68%  '
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('E10').Value = '  -3.63%  '
$ws.Range('E11').Value = '  -5.07%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -5.58%  '
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('E16').Value = '  -4.51%  '
$ws.Range('E17').Value = '  -1.77%  '
$ws.Range('E18').Value = '  +0.95%  '
$ws.Range('E19').Value = '  -3.22%  '
$ws.Range('E20').Value = '  +1.85%  '
$ws.Range('E21').Value = '  -3.86%  '
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('E23').Value = '  -3.86%  '
$ws.Range('E24').Value = '  -2.77%  '
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -3.59%  '
$ws.Range('E28').Value = '  -2.24%  '
$ws.Range('E29').Value = '  +8.10%  '
$ws.Range('E30').Value = '  -5.82%  '
$ws.Range('E31').Value = '  -9.92%  '
$ws.Range('E32').Value = '  -8.56%  '
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('E35').Value = '  -4.62%  '
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('E37').Value = '  -4.90%  '
$ws.Range('E38').Value = '  -6.24%  '
$ws.Range('E39').Value = '  -3.93%  '
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('E42').Value = '  -5.60%  '
$ws.Range('E43').Value = '  -5.11%  '
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('E45').Value = '  -3.82%  '
$ws.Range('E46').Value = '  -4.91%  '
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E49').Value = '  -3.96%  '
$ws.Range('E50').Value = '  -6.46%  '
$ws.Range('E51').Value = '  -3.95%  '

# Price column (D): values look numeric (and some use "." as a thousands
# separator mimicking the source site), so force text formatting first to
# avoid Excel auto-converting/rounding them, then restore default styling.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.154.88'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.643.46'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '307.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3899'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3865'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.000'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '49.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08689'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.74'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.126'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001293'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.464'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.644.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06902'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.907'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.160.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.328'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.791'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.42'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.495'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '140.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.353'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.414'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.823.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.974'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08052'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2685'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9523'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09206'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.03'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.464'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7562'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '13.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6918'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.466'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.089'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08410'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.265'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '133.52'
$ws.Range('D51').Style = 'Normal'
